$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 0.40000000000000002
$ws.Range("B2").Value = 32.753577351059569
$ws.Range("C2").Value = 16.290881626426543
$ws.Range("D2").Value = 0.49737717049400521
$ws.Range("E2").Value = 29.449268217519098
$ws.Range("F2").Value = 14.447644840319134
$ws.Range("G2").Value = 0.49059435818932723
$ws.Range("H2").Value = 299.08139733819786
$ws.Range("I2").Value = 267.0264231233715

# Row 3
$ws.Range("A3").Value = 0.44444444444444448
$ws.Range("B3").Value = 32.767804302535509
$ws.Range("C3").Value = 16.283798981148575
$ws.Range("D3").Value = 0.49694507544066863
$ws.Range("E3").Value = 29.436779645265183
$ws.Range("F3").Value = 14.538680507266456
$ws.Range("G3").Value = 0.49389507556425111
$ws.Range("H3").Value = 299
$ws.Range("I3").Value = 266.51806080714397

# Row 4
$ws.Range("A4").Value = 0.53333333333333344
$ws.Range("B4").Value = 32.790494709695217
$ws.Range("C4").Value = 16.275179168966535
$ws.Range("D4").Value = 0.49633832343964079
$ws.Range("E4").Value = 29.369643361785428
$ws.Range("F4").Value = 14.469798888897607
$ws.Range("G4").Value = 0.49267874010772339
$ws.Range("H4").Value = 299.08739985024806
$ws.Range("I4").Value = 266.63584998259716

# Row 5
$ws.Range("A5").Value = 0.57777777777777783
$ws.Range("B5").Value = 32.79853471903089
$ws.Range("C5").Value = 16.281372634773582
$ws.Range("D5").Value = 0.4964054880575669
$ws.Range("E5").Value = 29.349590373071486
$ws.Range("F5").Value = 14.46896269486612
$ws.Range("G5").Value = 0.49298687003623476
$ws.Range("H5").Value = 298.94543814924458
$ws.Range("I5").Value = 266.56821707743393

# Row 6
$ws.Range("A6").Value = 0.62222222222222223
$ws.Range("B6").Value = 32.805044014422748
$ws.Range("C6").Value = 16.287471593207371
$ws.Range("D6").Value = 0.4964929047510675
$ws.Range("E6").Value = 29.378356099740074
$ws.Range("F6").Value = 14.552087074166376
$ws.Range("G6").Value = 0.49533360630396628
$ws.Range("H6").Value = 299.05925983816877
$ws.Range("I6").Value = 265.92667667954333

# Row 7
$ws.Range("A7").Value = 0.66666666666666674
$ws.Range("B7").Value = 32.809679270121201
$ws.Range("C7").Value = 16.291820931320405
$ws.Range("D7").Value = 0.49655532433555005
$ws.Range("E7").Value = 29.36557716334837
$ws.Range("F7").Value = 14.55336341299337
$ws.Range("G7").Value = 0.49559262302385965
$ws.Range("H7").Value = 299.02033923187747
$ws.Range("I7").Value = 266.03366760463939

# Row 8
$ws.Range("A8").Value = 0.71111111111111114
$ws.Range("B8").Value = 32.762965081703925
$ws.Range("C8").Value = 16.248004781218622
$ws.Range("D8").Value = 0.49592595605127693
$ws.Range("E8").Value = 29.340203676791312
$ws.Range("F8").Value = 14.540525515770987
$ws.Range("G8").Value = 0.49558365974374041
$ws.Range("H8").Value = 299
$ws.Range("I8").Value = 265.54547235148067

# Row 9
$ws.Range("A9").Value = 0.75555555555555554
$ws.Range("B9").Value = 32.768806327497458
$ws.Range("C9").Value = 16.253001074109328
$ws.Range("D9").Value = 0.49599002513774398
$ws.Range("E9").Value = 29.318111259120577
$ws.Range("F9").Value = 14.587498380022573
$ws.Range("G9").Value = 0.49755928173867425
$ws.Range("H9").Value = 299
$ws.Range("I9").Value = 265.35560384384075

# Row 10
$ws.Range("A10").Value = 0.80000000000000004
$ws.Range("B10").Value = 32.772210078601958
$ws.Range("C10").Value = 16.25430715887039
$ws.Range("D10").Value = 0.49597836459260813
$ws.Range("E10").Value = 29.323472989536448
$ws.Range("F10").Value = 14.611649331637572
$ws.Range("G10").Value = 0.49829190890354208
$ws.Range("H10").Value = 298.90000000000003
$ws.Range("I10").Value = 265.58707292760351
